# se agrega la publicacion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing project name
$ws.Range("A2").Value = "TestRuben"

# Add the new published project row
$ws.Range("A3").Value = "TestRuben1"

# Reflect the new active selection left after the edit
$ws.Range("B12").Select()
